$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J, styled like the existing header (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-22
$data = @{
    2  = @(7, 8)
    3  = @(5, 6)
    4  = @(7, 7)
    5  = @(8, 8)
    6  = @(8, 8)
    7  = @(6, 6)
    8  = @(10, 10)
    9  = @(6, 6)
    10 = @(6, 6)
    11 = @(8, 8)
    12 = @(7, 7)
    13 = @(8, 8)
    14 = @(7, 7)
    15 = @(8, 8)
    16 = @(8, 8)
    17 = @(8, 8)
    18 = @(8, 8)
    19 = @(7, 8)
    20 = @(7, 8)
    21 = @(8, 8)
    22 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
